$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.842.77'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.636.07'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '215.27'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '0.5086'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.2582'
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('D9').Value = '0.06430'
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('D10').Value = '20.33'
$ws.Range('E10').Value = '  +4.65%  '
$ws.Range('D11').Value = '0.07800'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.251'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.641.26'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = '1.862.06'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '0.5593'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '0.0₅7670'
$ws.Range('E16').Value = '  +1.47%  '
$ws.Range('D17').Value = '63.27'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '25.851.31'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').Value = '4.374'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '192.51'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').Value = '9.936'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = '6.147'
$ws.Range('E23').Value = '  +2.39%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '1.758'
$ws.Range('E25').Value = '  -6.86%  '
$ws.Range('D26').Value = '139.00'
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').Value = '0.1235'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('D28').Value = '6.830'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').Value = '15.54'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '1.242'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = '0.04967'
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').Value = '3.306'
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').Value = '3.258'
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('D34').Value = '1.569'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = '0.9007'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('D37').Value = '0.5573'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D39').Value = '1.132.27'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').Value = '0.9963'
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').Value = '99.11'
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('D43').Value = '5.456'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('D44').Value = '0.8000'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('D46').Value = '55.65'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').Value = '0.4265'
$ws.Range('E47').Value = '  -3.80%  '
$ws.Range('D48').Value = '7.803'
$ws.Range('E48').Value = '  +3.28%  '
$ws.Range('D49').Value = '0.05029'
$ws.Range('E49').Value = '  -2.04%  '
$ws.Range('D50').Value = '0.9971'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('E51').Value = '  +0.38%  '
